$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A:D) -----------------------------------------------
# ColumnWidth is expressed in "characters"; the stored xml width is
# ColumnWidth + 5/7. Back the desired xml width out into the ColumnWidth
# value that needs to be supplied.
$ws.Columns.Item(1).ColumnWidth = (13.625 - 5.0/7.0)
$ws.Columns.Item(2).ColumnWidth = (18.5   - 5.0/7.0)
$ws.Columns.Item(3).ColumnWidth = (19.625 - 5.0/7.0)
$ws.Columns.Item(4).ColumnWidth = (16     - 5.0/7.0)

# --- Row heights (2:4) ---------------------------------------------------
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 18.75

# --- Publication-date column: swap the text dates for real dates --------
# and give each cell the banded blue/white look that Excel applied
# manually to this column.
$white = 16777215
$black = 0
$blue1 = 15189684   # FFB4C6E7
$blue2 = 15917529   # FFD9E1F2

function Format-DateCell($addr, $fillColor, $dateSerial, $dropBottomBorder) {
    $r = $ws.Range($addr)

    $r.Font.Size = 11
    $r.Font.Name = "游ゴシック"
    $r.Font.Color = $black

    $r.Interior.Color = $fillColor
    $r.Interior.PatternColor = $fillColor

    $r.Borders.LineStyle = 1
    $r.Borders.Color = $white
    if ($dropBottomBorder) {
        $r.Borders.Item(9).LineStyle = -4142
    }

    $r.NumberFormat = "mm-dd-yy"
    $r.Value = $dateSerial
}

Format-DateCell "C2" $blue1 44630 $false
Format-DateCell "C3" $blue2 44512 $false
Format-DateCell "C4" $blue1 44757 $true

# --- Selection -------------------------------------------------------------
$ws.Range("C2").Select()
